# daily auto push: 2026-01-14 18:49 UTC
# Insert two new daily rows (2026/01/14 23:00 and 2026/01/15 02:00) into the
# Sheet1 data table, pushing the existing rows 651-692 down to 653-694.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two full blank rows before the current row 651 (the 2026/12/29 block),
# right after the last existing 2026/01/14 entry (row 650). This shifts every
# row from 651 downward by two, matching the new dimension A1:D694.
$ws.Range("651:652").Insert()

# --- New row 651: 2026/01/14 (水), time 23:00, ranking 201 ---
$a651 = $ws.Range("A651")
$a651.NumberFormat = "@"
$a651.Value = "2026/01/14"
$a651.ClearFormats()
$ws.Range("B651").Value = "水"
$ws.Range("C651").Value = 23
$ws.Range("D651").Value = 201

# --- New row 652: 2026/01/15 (木), time 02:00, ranking 201 ---
$a652 = $ws.Range("A652")
$a652.NumberFormat = "@"
$a652.Value = "2026/01/15"
$a652.ClearFormats()
$ws.Range("B652").Value = "木"
$ws.Range("C652").Value = 2
$ws.Range("D652").Value = 201
